$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @(
    @("Column Enhanced Basic", 0.5089743589743591, 0.6923076923076923, 120.9329941272736),
    @("Column Enhanced BM25 Emphasis", 0.4717948717948718, 0.5846153846153846, 101.1998147964478),
    @("Column Enhanced Vector Emphasis", 0.4869230769230769, 0.6615384615384615, 70.17123675346375),
    @("Column Enhanced with Reranking", 0.71, 0.8615384615384616, 68.3456130027771)
)

$row = 11
foreach ($entry in $data) {
    $ws.Cells.Item($row, 1).Value = $entry[0]
    $ws.Cells.Item($row, 2).Value = $entry[1]
    $ws.Cells.Item($row, 3).Value = $entry[2]
    $ws.Cells.Item($row, 4).Value = $entry[3]
    $row++
}
